$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark. It marked the location of the
#    previous edit (end of the "UNDERGRAD - University of Western Ontario"
#    line). Word moves this bookmark to wherever text is typed next, so it
#    will be re-created at the new edit location further below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Fix the spelling error: "Certifcations" -> "Certifications" by
#    inserting a missing "i" in the middle of the word. We want the
#    resulting paragraph to end up as three separate runs (Certif / I /
#    cations) just like Word itself would leave behind after a live edit,
#    so we fence the freshly-typed character with a pair of temporary
#    bookmarks (added, then immediately removed) - this stops the engine
#    from silently re-coalescing the new run back into its identically
#    formatted neighbours.
$certRange = $d.Content
$certRange.Find.Execute("Certifcations", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insPos = $certRange.Start + "Certif".Length

$insRange = $d.Range($insPos, $insPos)
$insRange.InsertAfter("I")

$fenceA = $d.Range($insPos, $insPos)
$d.Bookmarks.Add("_TempFenceA", $fenceA)
$fenceB = $d.Range($insPos + 1, $insPos + 1)
$d.Bookmarks.Add("_TempFenceB", $fenceB)
$d.Bookmarks.Item("_TempFenceA").Delete()
$d.Bookmarks.Item("_TempFenceB").Delete()

# 3. Word also drops a "_GoBack" bookmark at the location of the most
#    recent edit. Re-create it between "com" and "mon" on the pivot-table
#    bullet point, matching where the edit actually happened. This splits
#    that run into two runs around the bookmark.
$bulletRange = $d.Content
$bulletRange.Find.Execute("Calculate margins and other common ratios using calculation on pivot table", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$markPos = $bulletRange.Start + "Calculate margins and other com".Length
$goBackRange = $d.Range($markPos, $markPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# Nudge the tail half of the split run (away from the bookmark boundary) so
# the engine treats it as freshly touched text, matching the rsid-less
# <w:r> Word itself leaves behind on that side of the split.
$touchPos = $markPos + 5
$touchIns = $d.Range($touchPos, $touchPos)
$touchIns.InsertAfter("X")
$touchDel = $d.Range($touchPos, $touchPos + 1)
$touchDel.Delete()
